$d = $word.ActiveDocument

# Replace the three plain "41 870 000" occurrences with "10 000 000"
$d.Content.Find.Execute("41 870 000", $true, $false, $false, $false, $false,
                         $true, 1, $false, "10 000 000", 2)

# Replace "-83 740 000" with "60 000 000"
$d.Content.Find.Execute("-83 740 000", $true, $false, $false, $false, $false,
                         $true, 1, $false, "60 000 000", 2)

# Replace "-41 870 000" with "70 000 000"
$d.Content.Find.Execute("-41 870 000", $true, $false, $false, $false, $false,
                         $true, 1, $false, "70 000 000", 2)

# Replace the Arabic spelled-out amount
$d.Content.Find.Execute("واحد وأربعون مليون وثمانمئة وسبعون ألف", $true, $false, $false, $false, $false,
                         $true, 1, $false, "عشرة مليون", 2)
